$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# --- Row 22: C22 numeric 4 -> text placeholder "0" (style matches D-column blank-data cells) ---
$ws.Range("C22").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("C22").PasteSpecial(-4122)

# --- Row 23: C23 text placeholder "0" -> numeric 1 ---
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("C23").Value = 1

# --- Row 30: F30 text placeholder "0" -> numeric 3 ---
$ws.Range("F30").NumberFormat = "#,##0"
$ws.Range("F30").Value = 3

# --- Bulk numeric updates ---
# Row 15
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 33.333333333333
$ws.Range("L15").Value = -14.285714285714
$ws.Range("M15").Value = 9.090909090909
$ws.Range("N15").Value = -65.714285714285

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 400
$ws.Range("F16").Value = 16
$ws.Range("H16").Value = 166.666666666667
$ws.Range("I16").Value = 129
$ws.Range("J16").Value = 148
$ws.Range("K16").Value = -12.837837837837
$ws.Range("L16").Value = -4.444444444444
$ws.Range("M16").Value = -30.27027027027
$ws.Range("N16").Value = -84.823529411764

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 172
$ws.Range("J17").Value = 177
$ws.Range("K17").Value = -2.824858757062
$ws.Range("L17").Value = -2.272727272727
$ws.Range("M17").Value = 72
$ws.Range("N17").Value = -61.607142857142

# Row 18
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 9
$ws.Range("H18").Value = -35.714285714285
$ws.Range("I18").Value = 163
$ws.Range("J18").Value = 173
$ws.Range("K18").Value = -5.780346820809
$ws.Range("L18").Value = 28.346456692913
$ws.Range("M18").Value = 32.520325203252
$ws.Range("N18").Value = -83.988212180746

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -53.846153846153
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -37.209302325581
$ws.Range("I19").Value = 485
$ws.Range("J19").Value = 505
$ws.Range("K19").Value = -3.960396039603
$ws.Range("L19").Value = 7.538802660753
$ws.Range("M19").Value = -5.088062622309
$ws.Range("N19").Value = -54.460093896713

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 133.333333333333
$ws.Range("I20").Value = 95
$ws.Range("J20").Value = 82
$ws.Range("K20").Value = 15.853658536585
$ws.Range("L20").Value = 13.095238095238
$ws.Range("M20").Value = 131.707317073171
$ws.Range("N20").Value = -89.947089947089

# Row 21
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -5
$ws.Range("F21").Value = 73
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = -14.117647058823
$ws.Range("I21").Value = 1058
$ws.Range("J21").Value = 1095
$ws.Range("K21").Value = -3.378995433789
$ws.Range("L21").Value = 7.085020242914
$ws.Range("M21").Value = 8.735868448098
$ws.Range("N21").Value = -75.855773619351

# Row 22
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 32
$ws.Range("K22").Value = -21.875
$ws.Range("M22").Value = -13.793103448275

# Row 23
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 107
$ws.Range("J23").Value = 119
$ws.Range("K23").Value = -10.084033613445
$ws.Range("L23").Value = 2.884615384615
$ws.Range("M23").Value = 44.594594594594

# Row 24
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 91
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = -26.612903225806
$ws.Range("I24").Value = 1472
$ws.Range("J24").Value = 1763
$ws.Range("K24").Value = -16.505955757232
$ws.Range("L24").Value = 21.152263374485
$ws.Range("M24").Value = 45.310957551826

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = 10.526315789473
$ws.Range("I25").Value = 286
$ws.Range("J25").Value = 288
$ws.Range("K25").Value = -0.694444444444
$ws.Range("L25").Value = 11.284046692607
$ws.Range("M25").Value = -8.626198083067

# Row 26
$ws.Range("C26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 200
$ws.Range("I26").Value = 21
$ws.Range("K26").Value = 40
$ws.Range("L26").Value = 16.666666666666

# Row 27
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 42
$ws.Range("K27").Value = -17.647058823529
$ws.Range("L27").Value = -19.230769230769

# Row 30
$ws.Range("I30").Value = 19
$ws.Range("K30").Value = 18.75
$ws.Range("L30").Value = 111.111111111111

